# fine schematico parte 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update quantity for "connettore 10 pin maschio" (D2): 3 -> 9
$ws.Range("D2").Value = 9

# Update the active cell/selection to D3, matching the saved view state
$ws.Activate()
$ws.Range("D3").Select()
